# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (G) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 1
$ws.Range("G3").Value  = 2
$ws.Range("G4").Value  = 3
$ws.Range("G5").Value  = 1
$ws.Range("G6").Value  = 2
$ws.Range("G7").Value  = 0
$ws.Range("G8").Value  = 5
$ws.Range("G9").Value  = 0
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
